$wb = $excel.ActiveWorkbook

# ===== Sheet ALC =====
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3253.4707
$ws.Range("J17").Value = 3253.4707
$ws.Range("L17").Value = 9760.4121
$ws.Range("N17").Value = -10096.4121
$ws.Range("H40").Value = 25018354
$ws.Range("I40").Value = 16942.715
$ws.Range("J40").Value = 83354984
$ws.Range("K40").Value = 16942.715
$ws.Range("L40").Value = 83354984
$ws.Range("M40").Value = -16767.715
$ws.Range("N40").Value = -83355334
$ws.Range("H64").Value = 25253332
$ws.Range("J64").Value = 50004296
$ws.Range("L64").Value = 50004296
$ws.Range("N64").Value = -50004792
$ws.Range("H67").Value = 25253332
$ws.Range("J67").Value = 50004296
$ws.Range("L67").Value = 50004296
$ws.Range("N67").Value = -50006012
$ws.Range("H76").Value = 4998.5
$ws.Range("J76").Value = 4998.5
$ws.Range("L76").Value = 4998.5
$ws.Range("N76").Value = -5628.5
$ws.Range("H79").Value = 4998.5
$ws.Range("J79").Value = 4998.5
$ws.Range("L79").Value = 4998.5
$ws.Range("N79").Value = -7182.5
$ws.Range("H106").Value = 111115110
$ws.Range("J106").Value = 10000
$ws.Range("L106").Value = 10000
$ws.Range("N106").Value = -11262
$ws.Range("H121").Value = 4700.5
$ws.Range("J121").Value = 4700.5
$ws.Range("L121").Value = 14101.5
$ws.Range("N121").Value = -17595.5

# ===== Sheet ARM =====
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 25284.945
$ws.Range("I32").Value = 28522.186
$ws.Range("K32").Value = 28522.186
$ws.Range("M32").Value = -28235.186
$ws.Range("H37").Value = 29999.8
$ws.Range("J37").Value = 29999.8
$ws.Range("L37").Value = 29999.8
$ws.Range("N37").Value = -30545.8
$ws.Range("H110").Value = 2553767.5
$ws.Range("I110").Value = 2915877.2
$ws.Range("J110").Value = 19000
$ws.Range("K110").Value = 2915877.2
$ws.Range("L110").Value = 19000
$ws.Range("M110").Value = -2913832.2
$ws.Range("N110").Value = -23090
$ws.Range("H132").Value = 12246.703
$ws.Range("J132").Value = 7604.273
$ws.Range("L132").Value = 22812.819
$ws.Range("N132").Value = -27872.819

# ===== Sheet BSM =====
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H110").Value = 88000
$ws.Range("J110").Value = 88000
$ws.Range("L110").Value = 88000
$ws.Range("N110").Value = -96180
$ws.Range("H134").Value = 3253.889
$ws.Range("I134").Value = 1612.1428
$ws.Range("K134").Value = 4836.428400000001
$ws.Range("M134").Value = -2301.428400000001

# ===== Sheet CRP =====
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 2231.6667
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 2231.6667
$ws.Range("K22").Value = 0
$ws.Range("M22").Value = 2231.6667
$ws.Range("N22").Value = -2931.6667
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("H41").Value = 300000
$ws.Range("I41").Value = 0
$ws.Range("K41").Value = 0

# ===== Sheet CUL =====
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1222.6666
$ws.Range("I5").Value = 786.2857
$ws.Range("J5").Value = 2750
$ws.Range("K5").Value = 2358.8571
$ws.Range("L5").Value = 8250
$ws.Range("M5").Value = -2246.8571
$ws.Range("N5").Value = -8474
$ws.Range("H56").Value = 333338140
$ws.Range("I56").Value = 333338140
$ws.Range("K56").Value = 333338140
$ws.Range("M56").Value = -333337610
$ws.Range("H74").Value = 21500.334
$ws.Range("I74").Value = 9501
$ws.Range("K74").Value = 28503
$ws.Range("M74").Value = -27442
$ws.Range("H77").Value = 21500.334
$ws.Range("I77").Value = 9501
$ws.Range("K77").Value = 85509
$ws.Range("M77").Value = -80205
$ws.Range("H86").Value = 370
$ws.Range("I86").Value = 322.5
$ws.Range("K86").Value = 967.5
$ws.Range("M86").Value = 218.5
$ws.Range("H89").Value = 370
$ws.Range("I89").Value = 322.5
$ws.Range("K89").Value = 2902.5
$ws.Range("M89").Value = 3025.5
$ws.Range("H122").Value = 1081.6666
$ws.Range("J122").Value = 1226.625
$ws.Range("L122").Value = 11039.625
$ws.Range("N122").Value = -15939.625
$ws.Range("H132").Value = 5883
$ws.Range("J132").Value = 13158
$ws.Range("L132").Value = 118422
$ws.Range("N132").Value = -123482
$ws.Range("H134").Value = 2940.8823
$ws.Range("I134").Value = 2940.8823
$ws.Range("K134").Value = 8822.6469
$ws.Range("M134").Value = -3752.6469
$ws.Range("H135").Value = 1222.6666
$ws.Range("I135").Value = 786.2857
$ws.Range("J135").Value = 2750
$ws.Range("K135").Value = 7076.571300000001
$ws.Range("L135").Value = 24750
$ws.Range("M135").Value = -4541.571300000001
$ws.Range("N135").Value = -29820

# ===== Sheet GSM =====
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 100000
$ws.Range("J47").Value = 100000
$ws.Range("L47").Value = 100000
$ws.Range("N47").Value = -101136
$ws.Range("H97").Value = 733.8461
$ws.Range("I97").Value = 614.1
$ws.Range("J97").Value = 1133
$ws.Range("K97").Value = 614.1
$ws.Range("L97").Value = 1133
$ws.Range("M97").Value = -118.1
$ws.Range("N97").Value = -2125

# ===== Sheet LTW =====
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 7409172
$ws.Range("I16").Value = 8334256
$ws.Range("K16").Value = 8334256
$ws.Range("M16").Value = -8334086
$ws.Range("H22").Value = 928.8333
$ws.Range("I22").Value = 814
$ws.Range("K22").Value = 814
$ws.Range("M22").Value = -519
$ws.Range("H27").Value = 928.8333
$ws.Range("I27").Value = 814
$ws.Range("K27").Value = 814
$ws.Range("M27").Value = -707
$ws.Range("H46").Value = 7172.533
$ws.Range("J46").Value = 7356.2856
$ws.Range("L46").Value = 7356.2856
$ws.Range("N46").Value = -7732.2856
$ws.Range("H61").Value = 4623.2
$ws.Range("I61").Value = 4096.875
$ws.Range("J61").Value = 5224.7144
$ws.Range("K61").Value = 4096.875
$ws.Range("L61").Value = 5224.7144
$ws.Range("M61").Value = -3894.875
$ws.Range("N61").Value = -5628.7144
$ws.Range("H68").Value = 3790339.2
$ws.Range("I68").Value = 4546805.5
$ws.Range("K68").Value = 4546805.5
$ws.Range("M68").Value = -4546056.5
$ws.Range("H71").Value = 3790339.2
$ws.Range("I71").Value = 4546805.5
$ws.Range("K71").Value = 22734027.5
$ws.Range("M71").Value = -22730283.5
$ws.Range("H82").Value = 7813587.5
$ws.Range("I82").Value = 15625450
$ws.Range("J82").Value = 1725
$ws.Range("K82").Value = 15625450
$ws.Range("L82").Value = 1725
$ws.Range("M82").Value = -15625089
$ws.Range("N82").Value = -2447
$ws.Range("H85").Value = 7813587.5
$ws.Range("I85").Value = 15625450
$ws.Range("J85").Value = 1725
$ws.Range("K85").Value = 15625450
$ws.Range("L85").Value = 1725
$ws.Range("M85").Value = -15624202
$ws.Range("N85").Value = -4221
$ws.Range("H113").Value = 4623.2
$ws.Range("I113").Value = 4096.875
$ws.Range("J113").Value = 5224.7144
$ws.Range("K113").Value = 4096.875
$ws.Range("L113").Value = 5224.7144
$ws.Range("M113").Value = -1926.875
$ws.Range("N113").Value = -9564.714400000001
$ws.Range("H132").Value = 4624.1035
$ws.Range("J132").Value = 4652.385
$ws.Range("L132").Value = 13957.155
$ws.Range("N132").Value = -19017.155

# ===== Sheet WVR =====
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H37").Value = 29990
$ws.Range("I37").Value = 29990
$ws.Range("K37").Value = 29990
$ws.Range("M37").Value = -29787
$ws.Range("H81").Value = 5216993
$ws.Range("J81").Value = 16947
$ws.Range("L81").Value = 33894
$ws.Range("N81").Value = -36016
$ws.Range("H84").Value = 5216993
$ws.Range("J84").Value = 16947
$ws.Range("L84").Value = 169470
$ws.Range("N84").Value = -180078
$ws.Range("H96").Value = 5743.727
$ws.Range("I96").Value = 5077
$ws.Range("J96").Value = 6299.3335
$ws.Range("K96").Value = 5077
$ws.Range("L96").Value = 6299.3335
$ws.Range("M96").Value = -3704
$ws.Range("N96").Value = -9045.333500000001
$ws.Range("H113").Value = 1119.4667
$ws.Range("I113").Value = 1224.625
$ws.Range("J113").Value = 999.2857
$ws.Range("K113").Value = 3673.875
$ws.Range("L113").Value = 2997.8571
$ws.Range("M113").Value = -1503.875
$ws.Range("N113").Value = -7337.8571
